$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# New wizard-related translation rows appended at the end of the Import sheet.
$rows = @(
    @{ Row = 187; Label = "lab.wizard.title"; Translation = "Průvodci" },
    @{ Row = 188; Label = "lab.wizard.subtitle"; Translation = "Veškeré užitečné postupy v aplikaci jsou řešené pomocí průvodců; tady je najdete." },
    @{ Row = 189; Label = "lab.wizard.build.title"; Translation = "Průvodce novým buildem" },
    @{ Row = 190; Label = "lab.wizard.build.subtitle"; Translation = "Tento průvodce vám pomůže zaevidovat nový build." }
)

foreach ($r in $rows) {
    $rowIndex = $r.Row
    $ws.Cells.Item($rowIndex, 1).Value = "cs"
    $ws.Cells.Item($rowIndex, 2).Value = $r.Label
    $ws.Cells.Item($rowIndex, 3).Value = $r.Translation

    $rng = $ws.Range("A" + $rowIndex + ":C" + $rowIndex)
    $rng.WrapText = $true
    $rng.Font.Size = 10
}

[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 175
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("B184").Select()
